$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testSheet")

# Move the value from C2 down to C4 (value previously read from C2,
# now read from C4 instead).
$ws.Range("C4").Value = $ws.Range("C2").Value2
$ws.Range("C2").ClearContents()

# Update the active selection to the new cell.
$ws.Range("C4").Select() | Out-Null
